$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 90
$ws.Range("I4").Value = 100
$ws.Range("J4").Value = 80
$ws.Range("K4").Value = 100
$ws.Range("L4").Value = 80
$ws.Range("M4").Value = 14
$ws.Range("N4").Value = -308
$ws.Range("H15").Value = 2452.9443
$ws.Range("I15").Value = 2452.9443
$ws.Range("K15").Value = 7358.8329
$ws.Range("M15").Value = -7189.8329
$ws.Range("H33").Value = 223.76471
$ws.Range("I33").Value = 232
$ws.Range("K33").Value = 232
$ws.Range("M33").Value = -3
$ws.Range("H98").Value = 874.125
$ws.Range("I98").Value = 874.125
$ws.Range("K98").Value = 874.125
$ws.Range("M98").Value = 623.875
$ws.Range("H107").Value = 1106.7858
$ws.Range("I107").Value = 1076.5385
$ws.Range("K107").Value = 1076.5385
$ws.Range("M107").Value = 843.4614999999999
$ws.Range("H116").Value = 5071.25
$ws.Range("I116").Value = 4994.5
$ws.Range("K116").Value = 4994.5
$ws.Range("M116").Value = -1552.5
$ws.Range("H122").Value = 874.125
$ws.Range("I122").Value = 874.125
$ws.Range("K122").Value = 2622.375
$ws.Range("M122").Value = -172.375
$ws.Range("H132").Value = 1041.963
$ws.Range("I132").Value = 1059.0834
$ws.Range("J132").Value = 905
$ws.Range("K132").Value = 3177.2502
$ws.Range("L132").Value = 2715
$ws.Range("M132").Value = -647.2502
$ws.Range("N132").Value = -7775
$ws.Range("H135").Value = 6788.3335
$ws.Range("I135").Value = 7866
$ws.Range("K135").Value = 70794
$ws.Range("M135").Value = -68259
$ws.Range("H137").Value = 3779.5
$ws.Range("I137").Value = 3324.25
$ws.Range("K137").Value = 9972.75
$ws.Range("M137").Value = -7422.75
$ws.Range("H141").Value = 4502.2856
$ws.Range("I141").Value = 3419.4167
$ws.Range("K141").Value = 10258.2501
$ws.Range("M141").Value = -5078.250100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3072.4075
$ws.Range("J32").Value = 2003
$ws.Range("L32").Value = 2003
$ws.Range("N32").Value = -2577
$ws.Range("H74").Value = 2494.5334
$ws.Range("I74").Value = 741.9
$ws.Range("K74").Value = 741.9
$ws.Range("M74").Value = 132.1
$ws.Range("H77").Value = 2494.5334
$ws.Range("I77").Value = 741.9
$ws.Range("K77").Value = 3709.5
$ws.Range("M77").Value = 658.5
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H132").Value = 2212
$ws.Range("I132").Value = 1782.3529
$ws.Range("J132").Value = 2876
$ws.Range("K132").Value = 5347.0587
$ws.Range("L132").Value = 8628
$ws.Range("M132").Value = -2817.0587
$ws.Range("N132").Value = -13688

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 4607.952
$ws.Range("I107").Value = 4629.8945
$ws.Range("J107").Value = 4399.5
$ws.Range("K107").Value = 4629.8945
$ws.Range("L107").Value = 4399.5
$ws.Range("M107").Value = -2709.8945
$ws.Range("N107").Value = -8239.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6952.1763
$ws.Range("I31").Value = 1784.1428
$ws.Range("J31").Value = 10569.8
$ws.Range("K31").Value = 1784.1428
$ws.Range("L31").Value = 10569.8
$ws.Range("M31").Value = -1489.1428
$ws.Range("N31").Value = -11159.8
$ws.Range("H34").Value = 6952.1763
$ws.Range("I34").Value = 1784.1428
$ws.Range("J34").Value = 10569.8
$ws.Range("K34").Value = 1784.1428
$ws.Range("L34").Value = 10569.8
$ws.Range("M34").Value = -1582.1428
$ws.Range("N34").Value = -10973.8
$ws.Range("H132").Value = 2681.75
$ws.Range("I132").Value = 1470.8334
$ws.Range("J132").Value = 4498.125
$ws.Range("K132").Value = 4412.5002
$ws.Range("L132").Value = 13494.375
$ws.Range("M132").Value = -1882.5002
$ws.Range("N132").Value = -18554.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 377
$ws.Range("I47").Value = 495
$ws.Range("J47").Value = 200
$ws.Range("K47").Value = 1485
$ws.Range("L47").Value = 600
$ws.Range("M47").Value = -1054
$ws.Range("N47").Value = -1462
$ws.Range("H68").Value = 1266.375
$ws.Range("J68").Value = 1275.8572
$ws.Range("L68").Value = 3827.5716
$ws.Range("N68").Value = -5449.571599999999
$ws.Range("H71").Value = 1266.375
$ws.Range("J71").Value = 1275.8572
$ws.Range("L71").Value = 11482.7148
$ws.Range("N71").Value = -19594.7148
$ws.Range("I86").Value = 500
$ws.Range("J86").Value = 600
$ws.Range("K86").Value = 1500
$ws.Range("L86").Value = 1800
$ws.Range("M86").Value = -314
$ws.Range("N86").Value = -4172
$ws.Range("I89").Value = 500
$ws.Range("J89").Value = 600
$ws.Range("K89").Value = 4500
$ws.Range("L89").Value = 5400
$ws.Range("M89").Value = 1428
$ws.Range("N89").Value = -17256
$ws.Range("H92").Value = 1900.2
$ws.Range("J92").Value = 3000
$ws.Range("L92").Value = 9000
$ws.Range("N92").Value = -11496
$ws.Range("H113").Value = 825.5455
$ws.Range("I113").Value = 569
$ws.Range("J113").Value = 1039.3334
$ws.Range("K113").Value = 1707
$ws.Range("L113").Value = 3118.0002
$ws.Range("M113").Value = 463
$ws.Range("N113").Value = -7458.0002
$ws.Range("H121").Value = 852.0625
$ws.Range("J121").Value = 921.36365
$ws.Range("L121").Value = 2764.09095
$ws.Range("N121").Value = -5384.09095
$ws.Range("H132").Value = 4541
$ws.Range("I132").Value = 2183
$ws.Range("J132").Value = 6899
$ws.Range("K132").Value = 19647
$ws.Range("L132").Value = 62091
$ws.Range("M132").Value = -17117
$ws.Range("N132").Value = -67151

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3945.625
$ws.Range("I102").Value = 3080.7144
$ws.Range("K102").Value = 3080.7144
$ws.Range("M102").Value = -1458.7144
$ws.Range("H122").Value = 1749.75
$ws.Range("I122").Value = 1749.75
$ws.Range("K122").Value = 5249.25
$ws.Range("M122").Value = -2799.25
$ws.Range("H126").Value = 2499.5
$ws.Range("I126").Value = 2499.5
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 7498.5
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -5028.5
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 3150.2122
$ws.Range("I132").Value = 2638.04
$ws.Range("J132").Value = 4750.75
$ws.Range("K132").Value = 7914.12
$ws.Range("L132").Value = 14252.25
$ws.Range("M132").Value = -5384.12
$ws.Range("N132").Value = -19312.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3000
$ws.Range("I7").Value = 3000
$ws.Range("K7").Value = 3000
$ws.Range("M7").Value = -2888
$ws.Range("H40").Value = 3162.353
$ws.Range("I40").Value = 3183.6428
$ws.Range("J40").Value = 3063
$ws.Range("K40").Value = 3183.6428
$ws.Range("L40").Value = 3063
$ws.Range("M40").Value = -3047.6428
$ws.Range("N40").Value = -3335
$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 3000
$ws.Range("K126").Value = 9000
$ws.Range("M126").Value = -6530
$ws.Range("H132").Value = 3304.4119
$ws.Range("I132").Value = 2475.5386
$ws.Range("K132").Value = 7426.6158
$ws.Range("M132").Value = -4896.6158

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1429.4117
$ws.Range("I122").Value = 1268.8125
$ws.Range("K122").Value = 3806.4375
$ws.Range("M122").Value = -1356.4375
$ws.Range("H126").Value = 1970.8
$ws.Range("I126").Value = 1970.8
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 5912.4
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -3442.4
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 2243.5134
$ws.Range("I132").Value = 1990.2142
$ws.Range("J132").Value = 3031.5557
$ws.Range("K132").Value = 5970.642599999999
$ws.Range("L132").Value = 9094.667099999999
$ws.Range("M132").Value = -3440.642599999999
$ws.Range("N132").Value = -14154.6671
